$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "is_locked" / "is_enabled" option columns: the remaining
# columns shift left. "order_by" moves from D1 to B1, "rem" moves from E1
# to C1, and a brand-new "tenant_id" column takes the D1 slot. The old E1
# cell no longer exists.
$ws.Range("B1").Value = '<%=comment.order_by%>'
$ws.Range("C1").Value = '<%=comment.rem%>'
$ws.Range("D1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'
$ws.Range("E1").ClearContents()
